{"js": "// Word JS API (Office.js) edit script.\n// The original document ends with an (almost) empty paragraph that only\n// contains a single tab run, immediately before the section break. The\n// author expanded that single paragraph into five paragraphs of new\n// personal-statement content (with first-line indents, a reused tab\n// run, and spell/grammar-check annotations preserved verbatim from the\n// authored OOXML).\n//\n// We locate that trailing paragraph and replace it (in place) with the\n// fully-authored OOXML for the new paragraphs via insertOoxml(..., Replace).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst targetParagraph = paragraphs.items[paragraphs.items.length - 1];\nconst targetRange = targetParagraph.getRange();\n\nconst newParagraphsOoxml = `<?xml version=\"1.0\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:spacing w:line=\"480\" w:lineRule=\"auto\"/>\n              <w:ind w:firstLine=\"720\"/>\n            </w:pPr>\n            <w:r>\n              <w:t xml:space=\"preserve\">In terms of working with a diverse group of students, I believe that diversity is essential for creating better solutions. Different perspectives and experiences can lead to new ideas and ways of thinking that can take a project or a problem in a new direction. I </w:t>\n            </w:r>\n            <w:r>\n              <w:lastRenderedPageBreak/>\n              <w:t>plan to work with a diverse group of students by being open to new ideas and perspectives, and by actively seeking out and valuing the contributions of all members of the group.</w:t>\n            </w:r>\n          </w:p>\n          <w:p>\n            <w:pPr>\n              <w:spacing w:line=\"480\" w:lineRule=\"auto\"/>\n            </w:pPr>\n          </w:p>\n          <w:p>\n            <w:pPr>\n              <w:spacing w:line=\"480\" w:lineRule=\"auto\"/>\n              <w:ind w:firstLine=\"720\"/>\n            </w:pPr>\n            <w:r>\n              <w:t>Personally, I bring a diverse perspective to the field of computer science. I come from a non-traditional background, and I have always been interested in how technology can be used to solve real-world problems. I am also fluent in multiple languages, which allows me to communicate and collaborate with people from different cultures and backgrounds.</w:t>\n            </w:r>\n            <w:r>\n              <w:tab/>\n            </w:r>\n          </w:p>\n          <w:p>\n            <w:pPr>\n              <w:spacing w:line=\"480\" w:lineRule=\"auto\"/>\n              <w:ind w:firstLine=\"720\"/>\n            </w:pPr>\n          </w:p>\n          <w:p>\n            <w:pPr>\n              <w:spacing w:line=\"480\" w:lineRule=\"auto\"/>\n              <w:ind w:firstLine=\"720\"/>\n            </w:pPr>\n            <w:r>\n              <w:t xml:space=\"preserve\">During the last winter </w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r>\n              <w:t>brea</w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r>\n              <w:t xml:space=\"preserve\">. I </w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r>\n              <w:t>surees</w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r>\n              <w:t xml:space=\"preserve\"> </w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r>\n              <w:t>deplove</w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r>\n              <w:t xml:space=\"preserve\"> my personal website. When I </w:t>\n            </w:r>\n            <w:proofErr w:type=\"gramStart\"/>\n            <w:r>\n              <w:t>doing</w:t>\n            </w:r>\n            <w:proofErr w:type=\"gramEnd\"/>\n            <w:r>\n              <w:t xml:space=\"preserve\"> this project. I </w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r>\n              <w:t>lerarn</w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r>\n              <w:t xml:space=\"preserve\"> a lot from the project</w:t>\n            </w:r>\n            <w:r>\n              <w:t xml:space=\"preserve\">. I </w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\ntargetRange.insertOoxml(newParagraphsOoxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Word COM interop edit script.\n# The original document ends with an (almost) empty paragraph that only\n# contains a single tab run, immediately before the section break. The\n# author expanded that single paragraph into five paragraphs of new\n# personal-statement content (with first-line indents, a reused tab\n# run, and spell/grammar-check annotations preserved verbatim from the\n# authored OOXML).\n#\n# We locate that trailing paragraph and replace it (in place) with the\n# fully-authored OOXML for the new paragraphs via Range.InsertXML(xml, \"Replace\"),\n# the COM analogue of Office.js's Range.insertOoxml(xml, Word.InsertLocation.replace).\n\n$d = $word.ActiveDocument\n\n$targetParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)\n$targetRange = $targetParagraph.Range\n\n$newParagraphsOoxml = @'\n<?xml version=\"1.0\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:spacing w:line=\"480\" w:lineRule=\"auto\"/>\n              <w:ind w:firstLine=\"720\"/>\n            </w:pPr>\n            <w:r>\n              <w:t xml:space=\"preserve\">In terms of working with a diverse group of students, I believe that diversity is essential for creating better solutions. Different perspectives and experiences can lead to new ideas and ways of thinking that can take a project or a problem in a new direction. I </w:t>\n            </w:r>\n            <w:r>\n              <w:lastRenderedPageBreak/>\n              <w:t>plan to work with a diverse group of students by being open to new ideas and perspectives, and by actively seeking out and valuing the contributions of all members of the group.</w:t>\n            </w:r>\n          </w:p>\n          <w:p>\n            <w:pPr>\n              <w:spacing w:line=\"480\" w:lineRule=\"auto\"/>\n            </w:pPr>\n          </w:p>\n          <w:p>\n            <w:pPr>\n              <w:spacing w:line=\"480\" w:lineRule=\"auto\"/>\n              <w:ind w:firstLine=\"720\"/>\n            </w:pPr>\n            <w:r>\n              <w:t>Personally, I bring a diverse perspective to the field of computer science. I come from a non-traditional background, and I have always been interested in how technology can be used to solve real-world problems. I am also fluent in multiple languages, which allows me to communicate and collaborate with people from different cultures and backgrounds.</w:t>\n            </w:r>\n            <w:r>\n              <w:tab/>\n            </w:r>\n          </w:p>\n          <w:p>\n            <w:pPr>\n              <w:spacing w:line=\"480\" w:lineRule=\"auto\"/>\n              <w:ind w:firstLine=\"720\"/>\n            </w:pPr>\n          </w:p>\n          <w:p>\n            <w:pPr>\n              <w:spacing w:line=\"480\" w:lineRule=\"auto\"/>\n              <w:ind w:firstLine=\"720\"/>\n            </w:pPr>\n            <w:r>\n              <w:t xml:space=\"preserve\">During the last winter </w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r>\n              <w:t>brea</w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r>\n              <w:t xml:space=\"preserve\">. I </w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r>\n              <w:t>surees</w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r>\n              <w:t xml:space=\"preserve\"> </w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r>\n              <w:t>deplove</w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r>\n              <w:t xml:space=\"preserve\"> my personal website. When I </w:t>\n            </w:r>\n            <w:proofErr w:type=\"gramStart\"/>\n            <w:r>\n              <w:t>doing</w:t>\n            </w:r>\n            <w:proofErr w:type=\"gramEnd\"/>\n            <w:r>\n              <w:t xml:space=\"preserve\"> this project. I </w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r>\n              <w:t>lerarn</w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r>\n              <w:t xml:space=\"preserve\"> a lot from the project</w:t>\n            </w:r>\n            <w:r>\n              <w:t xml:space=\"preserve\">. I </w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n'@\n\n$targetRange.InsertXML($newParagraphsOoxml, \"Replace\")\n"}
